$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-9 (columns C=runs, D=balls, E=fours, F=sixes)
# These reflect a re-ordering of the existing rows' stats.
$data = @{
    2 = @{ C = "13"; D = "8";  E = "1"; F = "1" }
    3 = @{ C = "9";  D = "11"; E = "1"; F = "0" }
    4 = @{ C = "5";  D = "3";  E = "1"; F = "0" }
    5 = @{ C = "2";  D = "4";  E = "0"; F = "0" }
    6 = @{ C = "11"; D = "11"; E = "2"; F = "0" }
    7 = @{ C = "12"; D = "9";  E = "1"; F = "1" }
    8 = @{ C = "16"; D = "10"; E = "2"; F = "1" }
    9 = @{ C = "25"; D = "11"; E = "1"; F = "3" }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("C$row").Value = "'" + $vals.C
    $ws.Range("D$row").Value = "'" + $vals.D
    $ws.Range("E$row").Value = "'" + $vals.E
    $ws.Range("F$row").Value = "'" + $vals.F
}
